# Fix setup of FCF (fund cash-flow) columns in the investor KYC import sheet.
# Adds two new columns - CF1 (a running sequence number) and CF2 (a letter
# code) - to the right of the existing "Agreement Committed Amount" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row
$ws.Range("S1").Value = "CF1"
$ws.Range("T1").Value = "CF2"

# Data rows
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = "A"

$ws.Range("S3").Value = 2
$ws.Range("T3").Value = "B"

$ws.Range("S4").Value = 3
$ws.Range("T4").Value = "C"

$ws.Range("S5").Value = 4
$ws.Range("T5").Value = "D"

# Scroll the view over to the newly added columns and move the selection
# to just below the last new column, matching the saved view state.
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("T6").Select()
